$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of columns M through Z (rows 1-10), which removes the
# extra indicator columns (Democratic Quality, Delivery Quality, etc.) while
# preserving the existing cell styles.
$ws.Range("M1:Z10").ClearContents()

# Update the active selection to F1 (single cell) instead of the whole
# column A1:A1048576 selection.
$ws.Range("F1").Select()
